# Auto-generated script to apply scheduled market-data refresh to Famfrit_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 7315.5557
$ws.Range("J9").Value = 3942.8572
$ws.Range("L9").Value = 3942.8572
$ws.Range("N9").Value = -4280.8572
$ws.Range("H19").Value = 1000
$ws.Range("J19").Value = 1000
$ws.Range("L19").Value = 1000
$ws.Range("N19").Value = -1350
$ws.Range("H32").Value = 4598.1113
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 4598.1113
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 4598.1113
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -5250.1113
$ws.Range("H38").Value = 4534.0557
$ws.Range("I38").Value = 2634.5
$ws.Range("J38").Value = 8333.166999999999
$ws.Range("K38").Value = 7903.5
$ws.Range("L38").Value = 24999.501
$ws.Range("M38").Value = -7531.5
$ws.Range("N38").Value = -25743.501
$ws.Range("H98").Value = 5583.6787
$ws.Range("I98").Value = 4542.3
$ws.Range("J98").Value = 8187.125
$ws.Range("K98").Value = 4542.3
$ws.Range("L98").Value = 8187.125
$ws.Range("M98").Value = -3044.3
$ws.Range("N98").Value = -11183.125
$ws.Range("H103").Value = 760.86664
$ws.Range("I103").Value = 1096.5
$ws.Range("J103").Value = 709.2308
$ws.Range("K103").Value = 3289.5
$ws.Range("L103").Value = 2127.6924
$ws.Range("M103").Value = -2703.5
$ws.Range("N103").Value = -3299.6924
$ws.Range("H116").Value = 8939.190000000001
$ws.Range("I116").Value = 8308.916999999999
$ws.Range("J116").Value = 9779.556
$ws.Range("K116").Value = 8308.916999999999
$ws.Range("L116").Value = 9779.556
$ws.Range("M116").Value = -4866.916999999999
$ws.Range("N116").Value = -16663.556
$ws.Range("H122").Value = 5583.6787
$ws.Range("I122").Value = 4542.3
$ws.Range("J122").Value = 8187.125
$ws.Range("K122").Value = 13626.9
$ws.Range("L122").Value = 24561.375
$ws.Range("M122").Value = -11176.9
$ws.Range("N122").Value = -29461.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2634.5454
$ws.Range("I45").Value = 1696.2
$ws.Range("K45").Value = 1696.2
$ws.Range("M45").Value = -1319.2
$ws.Range("H61").Value = 21740774
$ws.Range("I61").Value = 23811066
$ws.Range("K61").Value = 23811066
$ws.Range("M61").Value = -23810854
$ws.Range("H123").Value = 90429
$ws.Range("J123").Value = 90429
$ws.Range("L123").Value = 90429
$ws.Range("N123").Value = -100229
$ws.Range("H124").Value = 48551.75
$ws.Range("J124").Value = 48551.75
$ws.Range("L124").Value = 48551.75
$ws.Range("N124").Value = -58371.75
$ws.Range("H136").Value = 21740774
$ws.Range("I136").Value = 23811066
$ws.Range("K136").Value = 71433198
$ws.Range("M136").Value = -71430648

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2447.5264
$ws.Range("I20").Value = 1453.2727
$ws.Range("K20").Value = 1453.2727
$ws.Range("M20").Value = -1206.2727
$ws.Range("H94").Value = 1438.1852
$ws.Range("I94").Value = 1263.3334
$ws.Range("J94").Value = 1787.8889
$ws.Range("K94").Value = 1263.3334
$ws.Range("L94").Value = 1787.8889
$ws.Range("M94").Value = -812.3334
$ws.Range("N94").Value = -2689.8889
$ws.Range("H105").Value = 6191.973
$ws.Range("I105").Value = 8605.375
$ws.Range("K105").Value = 8605.375
$ws.Range("M105").Value = -6858.375
$ws.Range("H138").Value = 70780
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 70780
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 70780
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -81060
$ws.Range("H140").Value = 40709
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10874402
$ws.Range("I31").Value = 3291.3572
$ws.Range("K31").Value = 3291.3572
$ws.Range("M31").Value = -2996.3572
$ws.Range("H34").Value = 10874402
$ws.Range("I34").Value = 3291.3572
$ws.Range("K34").Value = 3291.3572
$ws.Range("M34").Value = -3089.3572
$ws.Range("H58").Value = 1223
$ws.Range("I58").Value = 1306.3636
$ws.Range("J58").Value = 1070.1666
$ws.Range("K58").Value = 1306.3636
$ws.Range("L58").Value = 1070.1666
$ws.Range("M58").Value = -1103.3636
$ws.Range("N58").Value = -1476.1666
$ws.Range("H86").Value = 2563.625
$ws.Range("I86").Value = 2334.8333
$ws.Range("K86").Value = 2334.8333
$ws.Range("M86").Value = -1211.8333
$ws.Range("H89").Value = 2563.625
$ws.Range("I89").Value = 2334.8333
$ws.Range("K89").Value = 11674.1665
$ws.Range("M89").Value = -6058.166499999999
$ws.Range("H105").Value = 26165
$ws.Range("I105").Value = 1549.6666
$ws.Range("K105").Value = 1549.6666
$ws.Range("M105").Value = 197.3334
$ws.Range("H132").Value = 226387.33
$ws.Range("I132").Value = 336743.16
$ws.Range("K132").Value = 1010229.48
$ws.Range("M132").Value = -1007699.48
$ws.Range("H136").Value = 1223
$ws.Range("I136").Value = 1306.3636
$ws.Range("J136").Value = 1070.1666
$ws.Range("K136").Value = 3919.0908
$ws.Range("L136").Value = 3210.4998
$ws.Range("M136").Value = -1369.0908
$ws.Range("N136").Value = -8310.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 97.92856999999999
$ws.Range("I2").Value = 75.90909000000001
$ws.Range("K2").Value = 455.4545400000001
$ws.Range("M2").Value = -342.4545400000001
$ws.Range("H44").Value = 6982.385
$ws.Range("I44").Value = 8485.684999999999
$ws.Range("J44").Value = 2902
$ws.Range("K44").Value = 25457.055
$ws.Range("L44").Value = 8706
$ws.Range("M44").Value = -25059.055
$ws.Range("N44").Value = -9502
$ws.Range("H88").Value = 20000
$ws.Range("J88").Value = 20000
$ws.Range("L88").Value = 60000
$ws.Range("N88").Value = -60856
$ws.Range("H91").Value = 20000
$ws.Range("J91").Value = 20000
$ws.Range("L91").Value = 60000
$ws.Range("N91").Value = -62964
$ws.Range("H97").Value = 1601.1
$ws.Range("J97").Value = 1388.375
$ws.Range("L97").Value = 4165.125
$ws.Range("N97").Value = -5157.125
$ws.Range("H140").Value = 2386
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6391.657
$ws.Range("I80").Value = 5705.381
$ws.Range("J80").Value = 7421.0713
$ws.Range("K80").Value = 5705.381
$ws.Range("L80").Value = 7421.0713
$ws.Range("M80").Value = -4707.381
$ws.Range("N80").Value = -9417.0713
$ws.Range("H83").Value = 6391.657
$ws.Range("I83").Value = 5705.381
$ws.Range("J83").Value = 7421.0713
$ws.Range("K83").Value = 28526.905
$ws.Range("L83").Value = 37105.35649999999
$ws.Range("M83").Value = -23534.905
$ws.Range("N83").Value = -47089.35649999999
$ws.Range("H132").Value = 3020.5
$ws.Range("I132").Value = 2689.4443
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 8068.3329
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -5538.3329
$ws.Range("N132").Value = -23060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 4136.273
$ws.Range("I20").Value = 2857
$ws.Range("J20").Value = 6375
$ws.Range("K20").Value = 2857
$ws.Range("L20").Value = 6375
$ws.Range("M20").Value = -2631
$ws.Range("N20").Value = -6827
$ws.Range("H22").Value = 815.125
$ws.Range("I22").Value = 953.4286
$ws.Range("J22").Value = 621.5
$ws.Range("K22").Value = 953.4286
$ws.Range("L22").Value = 621.5
$ws.Range("M22").Value = -658.4286
$ws.Range("N22").Value = -1211.5
$ws.Range("H27").Value = 815.125
$ws.Range("I27").Value = 953.4286
$ws.Range("J27").Value = 621.5
$ws.Range("K27").Value = 953.4286
$ws.Range("L27").Value = 621.5
$ws.Range("M27").Value = -846.4286
$ws.Range("N27").Value = -835.5
$ws.Range("H55").Value = 655.6429000000001
$ws.Range("I55").Value = 398.22223
$ws.Range("K55").Value = 398.22223
$ws.Range("M55").Value = -225.22223
$ws.Range("H61").Value = 3647.4707
$ws.Range("I61").Value = 2558.75
$ws.Range("K61").Value = 2558.75
$ws.Range("M61").Value = -2356.75
$ws.Range("H82").Value = 1798.2307
$ws.Range("J82").Value = 1236.6
$ws.Range("L82").Value = 1236.6
$ws.Range("N82").Value = -1958.6
$ws.Range("H85").Value = 1798.2307
$ws.Range("J85").Value = 1236.6
$ws.Range("L85").Value = 1236.6
$ws.Range("N85").Value = -3732.6
$ws.Range("H113").Value = 3647.4707
$ws.Range("I113").Value = 2558.75
$ws.Range("K113").Value = 2558.75
$ws.Range("M113").Value = -388.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4666.3335
$ws.Range("I81").Value = 2999.5
$ws.Range("K81").Value = 5999
$ws.Range("M81").Value = -4938
$ws.Range("H84").Value = 4666.3335
$ws.Range("I84").Value = 2999.5
$ws.Range("K84").Value = 29995
$ws.Range("M84").Value = -24691
$ws.Range("H132").Value = 3635.7856
$ws.Range("I132").Value = 2990.1
$ws.Range("K132").Value = 8970.299999999999
$ws.Range("M132").Value = -6440.299999999999
$ws.Range("H136").Value = 3999.5417
$ws.Range("I136").Value = 2686.875
$ws.Range("J136").Value = 6624.875
$ws.Range("K136").Value = 8060.625
$ws.Range("L136").Value = 19874.625
$ws.Range("M136").Value = -5510.625
$ws.Range("N136").Value = -24974.625

Write-Output "Applied scheduled market-data refresh"